$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "praneeth"
$ws.Range("B1").Value = "sai "
$ws.Range("C1").Value = "Javeed"

$ws.Range("C1").Select()
